$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide-number placeholder preview text: '<#>' -> '<Nr.>' (guillemets)
#    This token lives on the slide master, every slide layout, and the
#    notes master (wherever a "sldNum" placeholder is defined).
# ---------------------------------------------------------------------------
function Find-SlideNumberShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shape = $container.Shapes.Item($i)
        $isSlideNumPh = $false
        try {
            $isSlideNumPh = ($shape.PlaceholderFormat.Type -eq 13)   # ppPlaceholderSlideNumber
        } catch {
            $isSlideNumPh = $false
        }
        if ($isSlideNumPh) {
            return $shape
        }
    }
    return $null
}

$newSlideNumText = [char]0x2039 + "Nr." + [char]0x203A   # "<Nr.>"

# Slide master
$masterShape = Find-SlideNumberShape($p.SlideMaster)
if ($masterShape -ne $null) {
    $masterShape.TextFrame.TextRange.Text = $newSlideNumText
}

# Every slide layout attached to the master
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    $layoutShape = Find-SlideNumberShape($layout)
    if ($layoutShape -ne $null) {
        $layoutShape.TextFrame.TextRange.Text = $newSlideNumText
    }
}

# Notes master
$notesMasterShape = Find-SlideNumberShape($p.NotesMaster)
if ($notesMasterShape -ne $null) {
    $notesMasterShape.TextFrame.TextRange.Text = $newSlideNumText
}

# ---------------------------------------------------------------------------
# 2) Slide 12 ("Grafik 3" picture): reposition vertically (updated image).
#    off x stays 4087197 EMU, y moves from 2199788 -> 2129828 EMU.
# ---------------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
for ($i = 1; $i -le $slide12.Shapes.Count; $i++) {
    $shape = $slide12.Shapes.Item($i)
    if ($shape.Name -eq "Grafik 3") {
        $shape.Left = 4087197 / 12700
        $shape.Top = 2129828 / 12700
    }
}

Write-Host "Edit complete."
